$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "election": the single demo election becomes a real
# "Bundestagswahl" (type=2) instead of a generic placeholder.
# ---------------------------------------------------------------
$wsElection = $wb.Worksheets.Item("election")
$wsElection.Range("A2").Value = 2
$wsElection.Range("C2").Value = "Bundestagswahl"
$wsElection.Range("F2").Value = "2018-06-29T18:00:00.000Z"
$wsElection.Columns.Item(3).ColumnWidth = 14.666666666666668
$wsElection.Range("D18").Select()

# ---------------------------------------------------------------
# Sheet "candidate": add four new candidates (the 2017 Bundestag
# top candidates of AfD, CDU, SPD and FDP).
# ---------------------------------------------------------------
$wsCandidate = $wb.Worksheets.Item("candidate")

$wsCandidate.Range("A22").Value = "Gauland"
$wsCandidate.Range("B22").Value = "Alexander"
$wsCandidate.Range("C22").Value = "AFD"
$wsCandidate.Range("D22").Value = 2
$wsCandidate.Range("E22").Value = 0
$wsCandidate.Range("F22").Value = 2
$wsCandidate.Range("G22").Value = 60

$wsCandidate.Range("A23").Value = "Kramp-Karrenbauer"
$wsCandidate.Range("B23").Value = "Annegret"
$wsCandidate.Range("C23").Value = "CDU"
$wsCandidate.Range("D23").Value = 4
$wsCandidate.Range("E23").Value = 0
$wsCandidate.Range("F23").Value = 2
$wsCandidate.Range("G23").Value = 61

$wsCandidate.Range("A24").Value = "Nahles"
$wsCandidate.Range("B24").Value = "Andrea"
$wsCandidate.Range("C24").Value = "SPD"
$wsCandidate.Range("D24").Value = 4
$wsCandidate.Range("E24").Value = 0
$wsCandidate.Range("F24").Value = 2
$wsCandidate.Range("G24").Value = 62

$wsCandidate.Range("A25").Value = "Lindner"
$wsCandidate.Range("B25").Value = "Christian"
$wsCandidate.Range("C25").Value = "FDP"
$wsCandidate.Range("D25").Value = 4
$wsCandidate.Range("E25").Value = 0
$wsCandidate.Range("F25").Value = 2
$wsCandidate.Range("G25").Value = 63

$wsCandidate.Columns.Item(1).ColumnWidth = 17.666666666666664

# ---------------------------------------------------------------
# Sheet "party": add the four corresponding parties, and a
# duplicated "name" header in column G that mirrors column A.
# ---------------------------------------------------------------
$wsParty = $wb.Worksheets.Item("party")
$wsParty.Range("G1").Value = "name"

$wsParty.Range("B22").Value = "Alternative für Deutschland"
$wsParty.Range("C22").Value = 2
$wsParty.Range("D22").Value = 50
$wsParty.Range("E22").Value = 2
$wsParty.Range("F22").Value = 0
$wsParty.Range("G22").Value = "AFD"

$wsParty.Range("B23").Value = "Christlich Demokratische Union Deutschlands"
$wsParty.Range("C23").Value = 4
$wsParty.Range("D23").Value = 51
$wsParty.Range("E23").Value = 2
$wsParty.Range("F23").Value = 0
$wsParty.Range("G23").Value = "CDU"

$wsParty.Range("B24").Value = "Sozialdemokratische Partei Deutschlands"
$wsParty.Range("C24").Value = 4
$wsParty.Range("D24").Value = 52
$wsParty.Range("E24").Value = 2
$wsParty.Range("F24").Value = 0
$wsParty.Range("G24").Value = "SPD"

$wsParty.Range("B25").Value = "Freie Demokratische Partei"
$wsParty.Range("C25").Value = 4
$wsParty.Range("D25").Value = 53
$wsParty.Range("E25").Value = 2
$wsParty.Range("F25").Value = 0
$wsParty.Range("G25").Value = "FDP"

$wsParty.Range("C22:G25").Select()

# ---------------------------------------------------------------
# The "candidate" sheet becomes the active tab/selection.
# ---------------------------------------------------------------
$wsCandidate.Activate()
$wsCandidate.Range("C32").Select()
